# Goshop 13st Commit Rev 1.0
# Applies:
#  - Sheet1: replace row 21 (2025-03-09 -> 2025-03-12, new amounts) and
#    append row 22 with the same new values.
#  - Adds sheet "銷售記錄" with a per-file revenue breakdown.
#  - Adds sheet "銷售總合" with the total revenue.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 21 ---
# Force column A to stay plain text (matches existing date cells) instead
# of being auto-converted to a date serial number.
$ws1.Cells.Item(21, 1).NumberFormat = "@"
$ws1.Cells.Item(21, 1).Value = "2025-03-12"
$ws1.Cells.Item(21, 1).Style = "Normal"
$ws1.Cells.Item(21, 2).Value = 4024.81
$ws1.Cells.Item(21, 3).Value = 402.56
$ws1.Cells.Item(21, 4).Value = 3622.25

# --- Append new row 22 with the same data ---
$ws1.Cells.Item(22, 1).NumberFormat = "@"
$ws1.Cells.Item(22, 1).Value = "2025-03-12"
$ws1.Cells.Item(22, 1).Style = "Normal"
$ws1.Cells.Item(22, 2).Value = 4024.81
$ws1.Cells.Item(22, 3).Value = 402.56
$ws1.Cells.Item(22, 4).Value = 3622.25

# --- Add "銷售記錄" sheet (placed right after Sheet1) ---
$wsRecord = $wb.Worksheets.Add($null, $ws1)
$wsRecord.Name = "銷售記錄"

# Reuse the same header style as Sheet1's header row (bold, border, centered)
# by copying formats across instead of re-creating a new style.
$ws1.Range("A1:B1").Copy()
$wsRecord.Range("A1:B1").PasteSpecial(-4122)
$wsRecord.Cells.Item(1, 1).Value = "檔案名"
$wsRecord.Cells.Item(1, 2).Value = "revenue"

$records = @(
    @("goshop_orders_20250214.xlsx", 3563.31),
    @("goshop_orders_20250217.xlsx", 1471.14),
    @("goshop_orders_20250218_samuel-tw@outlook.com.xlsx", 1271.48),
    @("goshop_orders_20250219_samuel-tw@outlook.com.xlsx", 1292.45),
    @("goshop_orders_20250220_samuel-tw@outlook.com.xlsx", 1316.49),
    @("goshop_orders_20250221_samuel-tw@outlook.com.xlsx", 1539.97),
    @("goshop_orders_20250223_samuel-tw@outlook.com.xlsx", 1797.42),
    @("goshop_orders_20250225_samuel-tw@outlook.com.xlsx", 2084.96),
    @("goshop_orders_20250226_samuel-tw@outlook.com.xlsx", 1403.34),
    @("goshop_orders_20250227_samuel-tw@outlook.com.xlsx", 1903.14),
    @("goshop_orders_20250228_samuel-tw@outlook.com.xlsx", 1620.79),
    @("goshop_orders_20250301_samuel-tw@outlook.com.xlsx", 1663.16),
    @("goshop_orders_20250302_samuel-tw@outlook.com.xlsx", 1444.43),
    @("goshop_orders_20250303_samuel-tw@outlook.com.xlsx", 1603.7),
    @("goshop_orders_20250305_samuel-tw@outlook.com.xlsx", 3133.93),
    @("goshop_orders_20250306_samuel-tw@outlook.com.xlsx", 1407.38),
    @("goshop_orders_20250312_samuel-tw@outlook.com.xlsx", 3622.25)
)

$r = 2
foreach ($rec in $records) {
    $wsRecord.Cells.Item($r, 1).Value = $rec[0]
    $wsRecord.Cells.Item($r, 2).Value = $rec[1]
    $r = $r + 1
}

# --- Add "銷售總合" sheet (placed right after 銷售記錄) ---
$wsTotal = $wb.Worksheets.Add($null, $wsRecord)
$wsTotal.Name = "銷售總合"

$ws1.Range("A1").Copy()
$wsTotal.Range("A1").PasteSpecial(-4122)
$wsTotal.Cells.Item(1, 1).Value = "總收入"
$wsTotal.Cells.Item(2, 1).Value = 32139.34
